# Refresh the "Cryptos" price/volume table with the latest scrape.
# D column = Price (text, keeps thousand-dot formatting); E column = 1h volume change (text, "  +x.xx%  ").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.852.83"
$ws.Range("E2").Value = "  -0.91%  "
$ws.Range("D3").Value = "2.495.13"
$ws.Range("E3").Value = "  -1.20%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "534.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.76"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.18%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.996"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.24%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.565"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.38%  "
$ws.Range("D9").Value = "2.519.24"
$ws.Range("E9").Value = "  -0.50%  "
$ws.Range("E10").Value = "  +1.76%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.160"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.33%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.37"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.44%  "
$ws.Range("E13").Value = "  -2.32%  "
$ws.Range("D14").Value = "2.939.72"
$ws.Range("E14").Value = "  -1.10%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.19"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.26%  "
$ws.Range("D16").Value = "58.760.84"
$ws.Range("E16").Value = "  -0.95%  "
$ws.Range("E17").Value = "  -0.89%  "
$ws.Range("D18").Value = "2.507.89"
$ws.Range("E18").Value = "  -1.71%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.05"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.97%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.25"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.77%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "324.49"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.86%  "
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.86"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.15%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.70"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.45%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.419"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.60%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.165"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.38%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.995"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.24%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.59"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.77%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.76"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.28%  "
$ws.Range("D30").Value = "0.0₃0772"
$ws.Range("E30").Value = "  +0.26%  "
$ws.Range("E31").Value = "  -2.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "167.82"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.79%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.18"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.23%  "
$ws.Range("E34").Value = "  -0.14%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.40"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.86%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.56"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.14%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.10"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.44%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.57"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.36%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.81"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.59%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.826"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.62"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.72%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.27"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.30%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "279.89"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.70%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.995"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.29%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.605"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.97%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.86"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.04%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "128.19"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.54%  "
$ws.Range("E48").Value = "  +0.65%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0515"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.00%  "
$ws.Range("E50").Value = "  -0.38%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.36"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.01%  "
